$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared string rich-text runs) ---
# A8: "Volume 30   Number  37" -> "...38"  (change trailing run "37" -> "38")
$ws.Range("A8").Characters(21,2).Text = "38"
# C9: date range update
$ws.Range("C9").Characters(27,9).Text = "9/18/2023"
$ws.Range("C9").Characters(47,9).Text = "9/24/2023"

# --- Simple numeric updates (style/type unchanged) ---
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(16,4).Value = 4
$ws.Cells.Item(16,5).Value = -75
$ws.Cells.Item(16,6).Value = 12
$ws.Cells.Item(16,7).Value = 12
$ws.Cells.Item(16,8).Value = 0
$ws.Cells.Item(16,9).Value = 106
$ws.Cells.Item(16,10).Value = 139
$ws.Cells.Item(16,11).Value = -23.741007194244
$ws.Cells.Item(16,12).Value = 19.101123595505
$ws.Cells.Item(16,13).Value = -44.210526315789
$ws.Cells.Item(16,14).Value = -84.179104477611
$ws.Cells.Item(17,3).Value = 13
$ws.Cells.Item(17,5).Value = 160
$ws.Cells.Item(17,6).Value = 25
$ws.Cells.Item(17,7).Value = 17
$ws.Cells.Item(17,8).Value = 47.058823529411
$ws.Cells.Item(17,9).Value = 172
$ws.Cells.Item(17,10).Value = 140
$ws.Cells.Item(17,11).Value = 22.857142857142
$ws.Cells.Item(17,12).Value = 57.798165137614
$ws.Cells.Item(17,13).Value = 145.714285714286
$ws.Cells.Item(17,14).Value = -4.972375690607
$ws.Cells.Item(18,4).Value = 8
$ws.Cells.Item(18,5).Value = -75
$ws.Cells.Item(18,6).Value = 17
$ws.Cells.Item(18,7).Value = 26
$ws.Cells.Item(18,8).Value = -34.615384615384
$ws.Cells.Item(18,9).Value = 198
$ws.Cells.Item(18,10).Value = 167
$ws.Cells.Item(18,11).Value = 18.562874251497
$ws.Cells.Item(18,12).Value = 42.446043165467
$ws.Cells.Item(18,13).Value = -1.980198019801
$ws.Cells.Item(18,14).Value = -81.950774840474
$ws.Cells.Item(19,3).Value = 13
$ws.Cells.Item(19,4).Value = 14
$ws.Cells.Item(19,5).Value = -7.142857142857
$ws.Cells.Item(19,6).Value = 57
$ws.Cells.Item(19,7).Value = 61
$ws.Cells.Item(19,8).Value = -6.557377049180
$ws.Cells.Item(19,9).Value = 455
$ws.Cells.Item(19,10).Value = 490
$ws.Cells.Item(19,12).Value = 75.675675675675
$ws.Cells.Item(19,13).Value = 43.081761006289
$ws.Cells.Item(19,14).Value = 2.017937219730
$ws.Cells.Item(20,3).Value = 7
$ws.Cells.Item(20,4).Value = 2
$ws.Cells.Item(20,5).Value = 250
$ws.Cells.Item(20,6).Value = 23
$ws.Cells.Item(20,7).Value = 28
$ws.Cells.Item(20,8).Value = -17.857142857142
$ws.Cells.Item(20,9).Value = 270
$ws.Cells.Item(20,10).Value = 186
$ws.Cells.Item(20,11).Value = 45.161290322580
$ws.Cells.Item(20,12).Value = 143.243243243243
$ws.Cells.Item(20,13).Value = 66.666666666666
$ws.Cells.Item(20,14).Value = -92.512479201331
$ws.Cells.Item(21,3).Value = 36
$ws.Cells.Item(21,4).Value = 33
$ws.Cells.Item(21,5).Value = 9.090909090909
$ws.Cells.Item(21,6).Value = 135
$ws.Cells.Item(21,8).Value = -6.25
$ws.Cells.Item(21,9).Value = 1214
$ws.Cells.Item(21,10).Value = 1138
$ws.Cells.Item(21,11).Value = 6.678383128295
$ws.Cells.Item(21,12).Value = 67.448275862069
$ws.Cells.Item(21,13).Value = 27.253668763102
$ws.Cells.Item(21,14).Value = -79.870668214226
$ws.Cells.Item(22,6).Value = 4
$ws.Cells.Item(22,7).Value = 3
$ws.Cells.Item(22,8).Value = 33.333333333333
$ws.Cells.Item(22,9).Value = 33
$ws.Cells.Item(22,10).Value = 21
$ws.Cells.Item(22,11).Value = 57.142857142857
$ws.Cells.Item(22,12).Value = 10
$ws.Cells.Item(22,13).Value = 106.25
$ws.Cells.Item(23,6).Value = 5
$ws.Cells.Item(23,9).Value = 60
$ws.Cells.Item(23,11).Value = 50
$ws.Cells.Item(23,12).Value = 100
$ws.Cells.Item(23,13).Value = 130.769230769231
$ws.Cells.Item(24,3).Value = 20
$ws.Cells.Item(24,4).Value = 29
$ws.Cells.Item(24,5).Value = -31.034482758620
$ws.Cells.Item(24,6).Value = 121
$ws.Cells.Item(24,7).Value = 116
$ws.Cells.Item(24,8).Value = 4.310344827586
$ws.Cells.Item(24,9).Value = 1054
$ws.Cells.Item(24,10).Value = 1054
$ws.Cells.Item(24,11).Value = 0
$ws.Cells.Item(24,12).Value = 71.382113821138
$ws.Cells.Item(24,13).Value = 44.383561643835
$ws.Cells.Item(25,3).Value = 10
$ws.Cells.Item(25,4).Value = 5
$ws.Cells.Item(25,5).Value = 100
$ws.Cells.Item(25,6).Value = 36
$ws.Cells.Item(25,8).Value = 12.5
$ws.Cells.Item(25,9).Value = 342
$ws.Cells.Item(25,10).Value = 332
$ws.Cells.Item(25,11).Value = 3.012048192771
$ws.Cells.Item(25,12).Value = 40.163934426229
$ws.Cells.Item(25,13).Value = 16.326530612244
$ws.Cells.Item(26,6).Value = 3
$ws.Cells.Item(26,10).Value = 25
$ws.Cells.Item(26,11).Value = -12
$ws.Cells.Item(27,4).Value = 3
$ws.Cells.Item(27,5).Value = -100
$ws.Cells.Item(27,6).Value = 2
$ws.Cells.Item(27,7).Value = 8
$ws.Cells.Item(27,8).Value = -75
$ws.Cells.Item(27,10).Value = 48
$ws.Cells.Item(27,11).Value = -31.25
$ws.Cells.Item(27,12).Value = -13.157894736842
$ws.Cells.Item(28,8).Value = 0
$ws.Cells.Item(28,9).Value = 2
$ws.Cells.Item(28,11).Value = -77.777777777777
$ws.Cells.Item(28,12).Value = -50
$ws.Cells.Item(28,13).Value = -50
$ws.Cells.Item(28,14).Value = -80
$ws.Cells.Item(29,8).Value = 0
$ws.Cells.Item(29,9).Value = 2
$ws.Cells.Item(29,11).Value = -71.428571428571
$ws.Cells.Item(29,12).Value = -50
$ws.Cells.Item(29,13).Value = -50
$ws.Cells.Item(29,14).Value = -80

# --- Cells that change between numeric and text representation ---
# (these require a NumberFormat nudge + a format-only paste from a donor
#  cell that already has the exact target style, so the saved style index
#  matches what a human editing in Excel would have produced)
$ws.Cells.Item(15,3).NumberFormat = "@"
$ws.Cells.Item(15,3).Value = "0"
$ws.Range("C14").Copy()
$ws.Cells.Item(15,3).PasteSpecial(-4122)

$ws.Cells.Item(18,3).Value = 2
$ws.Range("G30").Copy()
$ws.Cells.Item(18,3).PasteSpecial(-4122)

$ws.Cells.Item(23,3).Value = 2
$ws.Range("G30").Copy()
$ws.Cells.Item(23,3).PasteSpecial(-4122)

$ws.Cells.Item(23,7).NumberFormat = "@"
$ws.Cells.Item(23,7).Value = "0"
$ws.Range("C14").Copy()
$ws.Cells.Item(23,7).PasteSpecial(-4122)

$ws.Cells.Item(23,8).NumberFormat = "@"
$ws.Cells.Item(23,8).Value = "***.*"
$ws.Range("C14").Copy()
$ws.Cells.Item(23,8).PasteSpecial(-4122)

$ws.Cells.Item(26,3).NumberFormat = "@"
$ws.Cells.Item(26,3).Value = "0"
$ws.Range("C14").Copy()
$ws.Cells.Item(26,3).PasteSpecial(-4122)

$ws.Cells.Item(26,4).Value = 1
$ws.Range("G30").Copy()
$ws.Cells.Item(26,4).PasteSpecial(-4122)

$ws.Cells.Item(26,5).Value = -100
$ws.Range("M14").Copy()
$ws.Cells.Item(26,5).PasteSpecial(-4122)

$ws.Cells.Item(26,7).Value = 1
$ws.Range("G30").Copy()
$ws.Cells.Item(26,7).PasteSpecial(-4122)

$ws.Cells.Item(26,8).Value = 200
$ws.Range("M14").Copy()
$ws.Cells.Item(26,8).PasteSpecial(-4122)

$ws.Cells.Item(27,3).NumberFormat = "@"
$ws.Cells.Item(27,3).Value = "0"
$ws.Range("C14").Copy()
$ws.Cells.Item(27,3).PasteSpecial(-4122)

$ws.Cells.Item(28,3).Value = 1
$ws.Range("G30").Copy()
$ws.Cells.Item(28,3).PasteSpecial(-4122)

$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = "0"
$ws.Range("C14").Copy()
$ws.Cells.Item(28,4).PasteSpecial(-4122)

$ws.Cells.Item(28,5).NumberFormat = "@"
$ws.Cells.Item(28,5).Value = "***.*"
$ws.Range("C14").Copy()
$ws.Cells.Item(28,5).PasteSpecial(-4122)

$ws.Cells.Item(28,6).Value = 1
$ws.Range("G30").Copy()
$ws.Cells.Item(28,6).PasteSpecial(-4122)

$ws.Cells.Item(29,3).Value = 1
$ws.Range("G30").Copy()
$ws.Cells.Item(29,3).PasteSpecial(-4122)

$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = "0"
$ws.Range("C14").Copy()
$ws.Cells.Item(29,4).PasteSpecial(-4122)

$ws.Cells.Item(29,5).NumberFormat = "@"
$ws.Cells.Item(29,5).Value = "***.*"
$ws.Range("C14").Copy()
$ws.Cells.Item(29,5).PasteSpecial(-4122)

$ws.Cells.Item(29,6).Value = 1
$ws.Range("G30").Copy()
$ws.Cells.Item(29,6).PasteSpecial(-4122)
